$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 6 new job-search log rows (135-140) below the current last row (134).
# A true row-insert on 135-138 inherits row 134s per-column formatting (date
# style on col A, hyperlink style on col F); the inherited hyperlink look on F
# is then stripped since the new Evidence cells are plain text in the source.
$ws.Range("A135:A138").EntireRow.Insert()
$ws.Range("F135:F138").ClearFormats()

# Row 135: TriCom Tech. Srvcs
$ws.Range("A135").Value = 45818
$ws.Range("B135").Value = "TriCom Tech. Srvcs"
$ws.Range("C135").Value = "Sr. Data Scientist"
$ws.Range("D135").Value = "Mpls, some kinda DS, want pyspark, etc. too"
$ws.Range("F135").Value = "https://tricomts.com/thanks?submissionGuid=b7e9bcc5-ecac-4bef-9e1c-185d1287974d"

# Row 136: Jacobs
$ws.Range("A136").Value = 45818
$ws.Range("B136").Value = "Jacobs"
$ws.Range("C136").Value = "Data Scientist"
$ws.Range("D136").Value = "Seattle, want 3 yrs sql and urban planning"
$ws.Range("F136").Value = "https://www.linkedin.com/jobs/view/4240985029/?trackingId=OFLdSA2TRi2%2Fs9YyFbSz9Q%3D%3D&refId=dvy6NrZQQQCjBMLe7xVjTw%3D%3D&midToken=AQFpnZsm4rTQjw&midSig=12PkMqoq3USrM1&trkEmail=eml-jobs_jymbii_digest-job_card-0-jobcard_body-null-1j75g~mbl0mm8z~1l-null-null&eid=1j75g-mbl0mm8z-1l&otpToken=MTAwMDE5ZTUxMTJhYzFjZWI1MjkwMWUzNDQxY2UwYjQ4N2M4ZDc0NDkxYTQ4NzZmNzdjMTAwNmU0NzViNTZmYWI2ODI5MGEzNDNkNWZlMGQ4Y2JjOTE2ZTg2M2NjNDIxNTRlMWI2YTdjOWI0Y2IsMSwx"

# Row 137: Versique
$ws.Range("A137").Value = 45818
$ws.Range("C137").Value = "Senior Data Scientist (ID:46984"
$ws.Range("B137").Value = "Versique"
$ws.Range("D137").Value = "Eden Prairie, MN"
$ws.Range("F137").Value = "https://www.linkedin.com/jobs/view/4242253047/?refId=ByteString(length%3D16%2Cbytes%3Dfabcb0fd...ab05cc36)&trackingId=SXyg3cJshTb0HVK2O8sqfQ%3D%3D"

# Row 138: Residio
$ws.Range("A138").Value = 45818
$ws.Range("B138").Value = "Residio"
$ws.Range("C138").Value = "Lead Data Scientist"
$ws.Range("D138").Value = "Eden Prairie, MN"
$ws.Range("F138").Value = "https://www.linkedin.com/jobs/view/4226498792/?trackingId=Ja9da%2FjH%2B28cifexcW8hSg%3D%3D&refId=ByteString%28length%3D16%2Cbytes%3Dfabcb0fd...ab05cc36%29&midToken=AQFpnZsm4rTQjw&midSig=1XrX5Uhf-9QrM1&trkEmail=eml-email_job_alert_digest_01-primary_job_list-0-jobcard_body-null-1j75g~mbivgu85~qt-null-null&eid=1j75g-mbivgu85-qt&otpToken=MTAwMDE5ZTUxMTJhYzFjZWI1MjkwMWUzNDYxN2UwYjE4YmNjZDU0MDk5YTQ4NzZmNzdjMTAwNmU0NzViNTY4YmFlZDFiNjgxNzJjNmQwMjBkM2E0YTQ3ZDEzMjMwMGNhZjUwYzlmNDZlOTA1ZWMsMSwx"

# Row 139: LTIMindtree
$ws.Range("B139").Value = "LTIMindtree"
$ws.Range("C139").Value = "Senior Data Scientist"
$ws.Range("D139").Value = "Bellevue, super low pay: `$110k"
$ws.Range("F139").Value = "https://www.linkedin.com/jobs/view/4240287783/?refId=ByteString(length%3D16%2Cbytes%3D8f2f56d1...a60e75b6)&trackingId=BYuILGb4CJcf1VeQWdf0sw%3D%3D"

# Row 140: Salesforce
$ws.Range("C140").Value = "Responsible AI Data Scientist - AI Red Teamer"
$ws.Range("B140").Value = "Salesforce"
$ws.Range("D140").Value = "Bellevue"
$ws.Range("F140").Value = "https://www.linkedin.com/jobs/view/4243029077/?trackingId=qdWducvyXbIrk%2FDsk13u4g%3D%3D&refId=ByteString%28length%3D16%2Cbytes%3D8f2f56d1...a60e75b6%29&midToken=AQFpnZsm4rTQjw&midSig=3VpEaHyaHgNrM1&trkEmail=eml-email_job_alert_digest_01-primary_job_list-0-jobcard_body-null-1j75g~mbfwat6u~q1-null-null&eid=1j75g-mbfwat6u-q1&otpToken=MTAwMDE5ZTUxMTJhYzFjZWI1MjkwMWUyNGYxZmUwYjA4OWNjZDY0MjlmYTQ4NzZmNzdjMTAwNmU0NzViNTZhMGE5YmRiNjlkNGJmOWU2NDNjY2QwM2FjMjY5M2MwYjIwMmFlMTk1M2Y4M2M5NWIsMSwx"

# Match the saved selection/cursor position from the edit.
$ws.Range("C141").Select() | Out-Null

